# Update the "built on" build-timestamp string that is embedded in several
# cells across the "About" and "Boundaries and methane sources" worksheets.
#
# Old value: January 30 2026 16.19.47 EST
# New value: February 02 2026 12.49.33 EST

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet -------------------------------------------------------

$a2 = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet ------------------------------
# Column S ("build_version") holds the same build string for every data row
# (rows 2 through 14).

for ($row = 2; $row -le 14; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # Column S = 19
    $val = $cell.Value()
    $cell.Value = $val.Replace($oldStamp, $newStamp)
}
